# Auto-generated Excel COM-interop script
# Updates market-price-derived columns (H-N) across all 8 sheets
# per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 65.40000000000001
$ws.Range("I8").Value = 65.40000000000001
$ws.Range("K8").Value = 196.2
$ws.Range("M8").Value = -57.20000000000002

$ws.Range("H38").Value = 4673.421
$ws.Range("J38").Value = 9484.166999999999
$ws.Range("L38").Value = 28452.501
$ws.Range("N38").Value = -29196.501

$ws.Range("H76").Value = 11558.625
$ws.Range("I76").Value = 11558.625
$ws.Range("K76").Value = 11558.625
$ws.Range("M76").Value = -11243.625

$ws.Range("H79").Value = 11558.625
$ws.Range("I79").Value = 11558.625
$ws.Range("K79").Value = 11558.625
$ws.Range("M79").Value = -10466.625

$ws.Range("H86").Value = 2700
$ws.Range("I86").Value = 2800.25
$ws.Range("J86").Value = 2566.3333
$ws.Range("K86").Value = 2800.25
$ws.Range("L86").Value = 2566.3333
$ws.Range("M86").Value = -1677.25
$ws.Range("N86").Value = -4812.3333

$ws.Range("H89").Value = 2700
$ws.Range("I89").Value = 2800.25
$ws.Range("J89").Value = 2566.3333
$ws.Range("K89").Value = 14001.25
$ws.Range("L89").Value = 12831.6665
$ws.Range("M89").Value = -8385.25
$ws.Range("N89").Value = -24063.6665

$ws.Range("H92").Value = 1098.64
$ws.Range("I92").Value = 972.5909
$ws.Range("K92").Value = 972.5909
$ws.Range("M92").Value = 275.4091

$ws.Range("H98").Value = 446.43243
$ws.Range("I98").Value = 371.6129
$ws.Range("J98").Value = 833
$ws.Range("K98").Value = 371.6129
$ws.Range("L98").Value = 833
$ws.Range("M98").Value = 1126.3871
$ws.Range("N98").Value = -3829

$ws.Range("H122").Value = 446.43243
$ws.Range("I122").Value = 371.6129
$ws.Range("J122").Value = 833
$ws.Range("K122").Value = 1114.8387
$ws.Range("L122").Value = 2499
$ws.Range("M122").Value = 1335.1613
$ws.Range("N122").Value = -7399

$ws.Range("H138").Value = 25643034
$ws.Range("I138").Value = 1230.7916
$ws.Range("J138").Value = 66669916
$ws.Range("K138").Value = 3692.3748
$ws.Range("L138").Value = 200009748
$ws.Range("M138").Value = 1447.6252
$ws.Range("N138").Value = -200020028

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1564.2858
$ws.Range("I45").Value = 1075
$ws.Range("J45").Value = 4500
$ws.Range("K45").Value = 1075
$ws.Range("L45").Value = 4500
$ws.Range("M45").Value = -698
$ws.Range("N45").Value = -5254

$ws.Range("H61").Value = 3044
$ws.Range("I61").Value = 2767.7036
$ws.Range("J61").Value = 4109.7144
$ws.Range("K61").Value = 2767.7036
$ws.Range("L61").Value = 4109.7144
$ws.Range("M61").Value = -2555.7036
$ws.Range("N61").Value = -4533.7144

$ws.Range("H74").Value = 60691.953
$ws.Range("I74").Value = 63627.35
$ws.Range("K74").Value = 63627.35
$ws.Range("M74").Value = -62753.35

$ws.Range("H77").Value = 60691.953
$ws.Range("I77").Value = 63627.35
$ws.Range("K77").Value = 318136.75
$ws.Range("M77").Value = -313768.75

$ws.Range("H97").Value = 2129.125
$ws.Range("I97").Value = 2301.5454
$ws.Range("K97").Value = 2301.5454
$ws.Range("M97").Value = -1805.5454

$ws.Range("H109").Value = 119999.5
$ws.Range("J109").Value = 119999.5
$ws.Range("L109").Value = 119999.5
$ws.Range("N109").Value = -122773.5

$ws.Range("H132").Value = 27979.945
$ws.Range("I132").Value = 2079.9795
$ws.Range("K132").Value = 6239.9385
$ws.Range("M132").Value = -3709.9385

$ws.Range("H136").Value = 3044
$ws.Range("I136").Value = 2767.7036
$ws.Range("J136").Value = 4109.7144
$ws.Range("K136").Value = 8303.110799999999
$ws.Range("L136").Value = 12329.1432
$ws.Range("M136").Value = -5753.110799999999
$ws.Range("N136").Value = -17429.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1004
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H94").Value = 2602.611
$ws.Range("I94").Value = 2444.375
$ws.Range("K94").Value = 2444.375
$ws.Range("M94").Value = -1993.375

$ws.Range("H134").Value = 1361.1428
$ws.Range("I134").Value = 1390.2727
$ws.Range("J134").Value = 1254.3334
$ws.Range("K134").Value = 4170.8181
$ws.Range("L134").Value = 3763.0002
$ws.Range("M134").Value = -1635.8181
$ws.Range("N134").Value = -8833.0002

$ws.Range("H139").Value = 167499.5
$ws.Range("J139").Value = 167499.5
$ws.Range("L139").Value = 167499.5
$ws.Range("N139").Value = -177779.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 410
$ws.Range("I10").Value = 262.5
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 262.5
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = -123.5
$ws.Range("N10").Value = -1278

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H134").Value = 2215.647
$ws.Range("I134").Value = 1843.3077
$ws.Range("K134").Value = 5529.9231
$ws.Range("M134").Value = -2994.9231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 868.6
$ws.Range("J23").Value = 868.6
$ws.Range("L23").Value = 2605.8
$ws.Range("N23").Value = -3075.8

$ws.Range("H47").Value = 1180.4166
$ws.Range("I47").Value = 613.5
$ws.Range("J47").Value = 1747.3334
$ws.Range("K47").Value = 1840.5
$ws.Range("L47").Value = 5242.0002
$ws.Range("M47").Value = -1409.5
$ws.Range("N47").Value = -6104.0002

$ws.Range("H80").Value = 14144.833
$ws.Range("I80").Value = 3140.3333
$ws.Range("J80").Value = 25149.334
$ws.Range("K80").Value = 9420.999899999999
$ws.Range("L80").Value = 75448.00199999999
$ws.Range("M80").Value = -8484.999899999999
$ws.Range("N80").Value = -77320.00199999999

$ws.Range("H83").Value = 14144.833
$ws.Range("I83").Value = 3140.3333
$ws.Range("J83").Value = 25149.334
$ws.Range("K83").Value = 28262.9997
$ws.Range("L83").Value = 226344.006
$ws.Range("M83").Value = -23582.9997
$ws.Range("N83").Value = -235704.006

$ws.Range("H94").Value = 10001
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws.Range("H134").Value = 3729
$ws.Range("J134").Value = 9022
$ws.Range("L134").Value = 27066
$ws.Range("N134").Value = -37206

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H113").Value = 4014.4
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 4143
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 4143
$ws.Range("M113").Value = -1330
$ws.Range("N113").Value = -8483

$ws.Range("H123").Value = 37942.75
$ws.Range("I123").Value = 53333.332
$ws.Range("J123").Value = 32812.555
$ws.Range("K123").Value = 53333.332
$ws.Range("L123").Value = 32812.555
$ws.Range("M123").Value = -50883.332
$ws.Range("N123").Value = -37712.555

$ws.Range("H132").Value = 1299.1818
$ws.Range("I132").Value = 1389
$ws.Range("J132").Value = 895
$ws.Range("K132").Value = 4167
$ws.Range("L132").Value = 2685
$ws.Range("M132").Value = -1637
$ws.Range("N132").Value = -7745

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3665.889
$ws.Range("I93").Value = 3124.125
$ws.Range("K93").Value = 3124.125
$ws.Range("M93").Value = -1876.125

$ws.Range("H100").Value = 3578.7144
$ws.Range("I100").Value = 3300.5
$ws.Range("J100").Value = 4274.25
$ws.Range("K100").Value = 3300.5
$ws.Range("L100").Value = 4274.25
$ws.Range("M100").Value = -2759.5
$ws.Range("N100").Value = -5356.25

$ws.Range("H132").Value = 12441.889
$ws.Range("I132").Value = 3068.2856
$ws.Range("K132").Value = 9204.856800000001
$ws.Range("M132").Value = -6674.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 106
$ws.Range("I14").Value = 107.166664
$ws.Range("K14").Value = 107.166664
$ws.Range("M14").Value = 60.833336

$ws.Range("H101").Value = 28500
$ws.Range("J101").Value = 28500
$ws.Range("L101").Value = 28500
$ws.Range("N101").Value = -34990

$ws.Range("H132").Value = 1336.762
$ws.Range("I132").Value = 1130.6154
$ws.Range("J132").Value = 1671.75
$ws.Range("K132").Value = 3391.8462
$ws.Range("L132").Value = 5015.25
$ws.Range("M132").Value = -861.8462
$ws.Range("N132").Value = -10075.25
